# Helper: write a value into a cell while forcing TEXT storage (keeps
# leading zeros / numeric-looking strings like "005052" or "92.21" as text,
# matching the source workbook's inlineStr cells) and then clears the
# NumberFormat override so no stray style lingers on the cell.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet, positioned right after "总计" and
#    before "2022-Q2" (mirrors the new sheetId ordering in the diff).
#    Copying "2022-Q2" gives us an identical header row + styling
#    (s="2" borders/bold on row 1 and column A) for free.
# ------------------------------------------------------------------
$sourceQ2 = $wb.Worksheets.Item("2022-Q2")
$sourceQ2.Copy($sourceQ2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The copied sheet only has 2 data rows (rows 2-3); add a 3rd by cloning
# the formatting of row 3 down into row 4.
$q3.Range("A3:H3").Copy()
$q3.Range("A4:H4").PasteSpecial(-4122)

$q3.Range("A2").Value = 0
Set-TextCell $q3.Range("B2") "501305"
$q3.Range("C2").Value = "汇添富中证港股通高股息投资指数（LOF）A"
Set-TextCell $q3.Range("D2") "0.87"
Set-TextCell $q3.Range("E2") "92.21"
Set-TextCell $q3.Range("F2") "4.27"
Set-TextCell $q3.Range("G2") "0.0371"
$q3.Range("H2").Value = 4

$q3.Range("A3").Value = 1
Set-TextCell $q3.Range("B3") "513530"
$q3.Range("C3").Value = "华泰柏瑞中证港股通高股息投资ETF（QDII）"
Set-TextCell $q3.Range("D3") "0.78"
Set-TextCell $q3.Range("E3") "95.80"
Set-TextCell $q3.Range("F3") "4.45"
Set-TextCell $q3.Range("G3") "0.0347"
$q3.Range("H3").Value = 4

$q3.Range("A4").Value = 2
Set-TextCell $q3.Range("B4") "501306"
$q3.Range("C4").Value = "汇添富中证港股通高股息投资指数（LOF）C"
Set-TextCell $q3.Range("D4") "0.17"
Set-TextCell $q3.Range("E4") "92.21"
Set-TextCell $q3.Range("F4") "4.27"
Set-TextCell $q3.Range("G4") "0.0073"
$q3.Range("H4").Value = 4

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: add a 2022-Q3 row right after
#    the header and shift every existing row down by one, recomputing
#    the 0-based index in column A. Final table (rows 2-9):
#      0 2022-Q3 3 0.08
#      1 2022-Q2 2 0.12
#      2 2022-Q1 2 0.28
#      3 2021-Q4 4 0.04
#      4 2021-Q3 3 0.10
#      5 2021-Q2 1 0.01
#      6 2021-Q1 3 0.12
#      7 2020-Q4 2 0.09
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 9 doesn't exist yet - create it by cloning row 8's formatting
# first so column A keeps its s="2" style.
$total.Range("A8:D8").Copy()
$total.Range("A9:D9").PasteSpecial(-4122)

$summaryRows = @(
    @{B = "2022-Q3"; C = 3; D = 0.08},
    @{B = "2022-Q2"; C = 2; D = 0.12},
    @{B = "2022-Q1"; C = 2; D = 0.28},
    @{B = "2021-Q4"; C = 4; D = 0.04},
    @{B = "2021-Q3"; C = 3; D = 0.1},
    @{B = "2021-Q2"; C = 1; D = 0.01},
    @{B = "2021-Q1"; C = 3; D = 0.12},
    @{B = "2020-Q4"; C = 2; D = 0.09}
)

$r = 2
$idx = 0
foreach ($row in $summaryRows) {
    $total.Range("A$r").Value = $idx
    $total.Range("B$r").Value = $row.B
    $total.Range("C$r").Value = $row.C
    $total.Range("D$r").Value = $row.D
    $r = $r + 1
    $idx = $idx + 1
}

# Keep the originally-active tab ("2020-Q4") active - sheet creation
# (Add/Copy) otherwise steals activation for the newly inserted sheet.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "2022-Q3 sheet added"
